$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(81, 8).Value = 30000
$ws.Cells.Item(81, 10).Value = 30000
$ws.Cells.Item(81, 12).Value = 30000
$ws.Cells.Item(81, 14).Value = -31996
$ws.Cells.Item(84, 8).Value = 30000
$ws.Cells.Item(84, 10).Value = 30000
$ws.Cells.Item(84, 12).Value = 90000
$ws.Cells.Item(84, 14).Value = -99984
$ws.Cells.Item(123, 8).Value = 33126.785
$ws.Cells.Item(123, 10).Value = 33126.785
$ws.Cells.Item(123, 12).Value = 33126.785
$ws.Cells.Item(123, 14).Value = -42926.785
$ws.Cells.Item(124, 8).Value = 29666.666
$ws.Cells.Item(124, 10).Value = 29666.666
$ws.Cells.Item(124, 12).Value = 29666.666
$ws.Cells.Item(124, 14).Value = -39486.666
$ws.Cells.Item(126, 8).Value = 41306.668
$ws.Cells.Item(126, 10).Value = 41306.668
$ws.Cells.Item(126, 12).Value = 41306.668
$ws.Cells.Item(126, 14).Value = -51186.668
$ws.Cells.Item(132, 8).Value = 6900433
$ws.Cells.Item(132, 9).Value = 7410854
$ws.Cells.Item(132, 11).Value = 22232562
$ws.Cells.Item(132, 13).Value = -22230032
$ws.Cells.Item(133, 8).Value = 16866.674
$ws.Cells.Item(133, 10).Value = 16866.674
$ws.Cells.Item(133, 12).Value = 16866.674
$ws.Cells.Item(133, 14).Value = -26986.674
$ws.Cells.Item(134, 8).Value = 25534.482
$ws.Cells.Item(134, 10).Value = 25534.482
$ws.Cells.Item(134, 12).Value = 25534.482
$ws.Cells.Item(134, 14).Value = -35674.482
$ws.Cells.Item(136, 8).Value = 18458.586
$ws.Cells.Item(136, 10).Value = 18458.586
$ws.Cells.Item(136, 12).Value = 18458.586
$ws.Cells.Item(136, 14).Value = -28658.586
$ws.Cells.Item(137, 8).Value = 6212.136
$ws.Cells.Item(137, 9).Value = 9650.799999999999
$ws.Cells.Item(137, 10).Value = 3346.5833
$ws.Cells.Item(137, 11).Value = 28952.4
$ws.Cells.Item(137, 12).Value = 10039.7499
$ws.Cells.Item(137, 13).Value = -26402.4
$ws.Cells.Item(137, 14).Value = -15139.7499
$ws.Cells.Item(138, 8).Value = 3561.621
$ws.Cells.Item(138, 9).Value = 1905.5918
$ws.Cells.Item(138, 10).Value = 5325.6523
$ws.Cells.Item(138, 11).Value = 5716.7754
$ws.Cells.Item(138, 12).Value = 15976.9569
$ws.Cells.Item(138, 13).Value = -576.7753999999995
$ws.Cells.Item(138, 14).Value = -26256.9569
$ws.Cells.Item(139, 8).Value = 19722.908
$ws.Cells.Item(139, 10).Value = 19722.908
$ws.Cells.Item(139, 12).Value = 19722.908
$ws.Cells.Item(139, 14).Value = -30002.908
$ws.Cells.Item(141, 8).Value = 607722.4
$ws.Cells.Item(141, 9).Value = 1304.3529
$ws.Cells.Item(141, 11).Value = 3913.0587
$ws.Cells.Item(141, 13).Value = 1266.9413

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15249.913
$ws.Cells.Item(32, 9).Value = 14779.786
$ws.Cells.Item(32, 10).Value = 15981.223
$ws.Cells.Item(32, 11).Value = 14779.786
$ws.Cells.Item(32, 12).Value = 15981.223
$ws.Cells.Item(32, 13).Value = -14492.786
$ws.Cells.Item(32, 14).Value = -16555.223
$ws.Cells.Item(135, 8).Value = 21220.334
$ws.Cells.Item(135, 10).Value = 21220.334
$ws.Cells.Item(135, 12).Value = 21220.334
$ws.Cells.Item(135, 14).Value = -31360.334
$ws.Cells.Item(137, 8).Value = 17688.3
$ws.Cells.Item(137, 10).Value = 17688.3
$ws.Cells.Item(137, 12).Value = 17688.3
$ws.Cells.Item(137, 14).Value = -27888.3
$ws.Cells.Item(139, 8).Value = 16275.279
$ws.Cells.Item(139, 10).Value = 16275.279
$ws.Cells.Item(139, 12).Value = 16275.279
$ws.Cells.Item(139, 14).Value = -26555.279

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1461.5264
$ws.Cells.Item(105, 9).Value = 1359.875
$ws.Cells.Item(105, 10).Value = 2003.6666
$ws.Cells.Item(105, 11).Value = 1359.875
$ws.Cells.Item(105, 12).Value = 2003.6666
$ws.Cells.Item(105, 13).Value = 387.125
$ws.Cells.Item(105, 14).Value = -5497.6666
$ws.Cells.Item(122, 8).Value = 29666.666
$ws.Cells.Item(122, 10).Value = 29666.666
$ws.Cells.Item(122, 12).Value = 29666.666
$ws.Cells.Item(122, 14).Value = -39466.666
$ws.Cells.Item(126, 8).Value = 25000
$ws.Cells.Item(126, 10).Value = 25000
$ws.Cells.Item(126, 12).Value = 25000
$ws.Cells.Item(126, 14).Value = -34880
$ws.Cells.Item(130, 8).Value = 29980
$ws.Cells.Item(130, 10).Value = 29980
$ws.Cells.Item(130, 12).Value = 29980
$ws.Cells.Item(130, 14).Value = -40020
$ws.Cells.Item(132, 8).Value = 21286.285
$ws.Cells.Item(132, 10).Value = 21286.285
$ws.Cells.Item(132, 12).Value = 21286.285
$ws.Cells.Item(132, 14).Value = -31406.285
$ws.Cells.Item(135, 8).Value = 18023.902
$ws.Cells.Item(135, 10).Value = 17840.967
$ws.Cells.Item(135, 12).Value = 17840.967
$ws.Cells.Item(135, 14).Value = -27980.967
$ws.Cells.Item(138, 8).Value = 15000.837
$ws.Cells.Item(138, 10).Value = 15000.837
$ws.Cells.Item(138, 12).Value = 15000.837
$ws.Cells.Item(138, 14).Value = -25280.837
$ws.Cells.Item(140, 8).Value = 18875.465
$ws.Cells.Item(140, 10).Value = 18875.465
$ws.Cells.Item(140, 12).Value = 18875.465
$ws.Cells.Item(140, 14).Value = -29235.465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1026.625
$ws.Cells.Item(22, 9).Value = 343.83334
$ws.Cells.Item(22, 11).Value = 343.83334
$ws.Cells.Item(22, 13).Value = 6.166659999999979
$ws.Cells.Item(31, 8).Value = 4424
$ws.Cells.Item(31, 9).Value = 3485.9048
$ws.Cells.Item(31, 10).Value = 5059.484
$ws.Cells.Item(31, 11).Value = 3485.9048
$ws.Cells.Item(31, 12).Value = 5059.484
$ws.Cells.Item(31, 13).Value = -3190.9048
$ws.Cells.Item(31, 14).Value = -5649.484
$ws.Cells.Item(34, 8).Value = 4424
$ws.Cells.Item(34, 9).Value = 3485.9048
$ws.Cells.Item(34, 10).Value = 5059.484
$ws.Cells.Item(34, 11).Value = 3485.9048
$ws.Cells.Item(34, 12).Value = 5059.484
$ws.Cells.Item(34, 13).Value = -3283.9048
$ws.Cells.Item(34, 14).Value = -5463.484
$ws.Cells.Item(94, 8).Value = 6180.9414
$ws.Cells.Item(94, 10).Value = 5396.0835
$ws.Cells.Item(94, 12).Value = 5396.0835
$ws.Cells.Item(94, 14).Value = -6298.0835
$ws.Cells.Item(130, 8).Value = 30000
$ws.Cells.Item(130, 10).Value = 30000
$ws.Cells.Item(130, 12).Value = 30000
$ws.Cells.Item(130, 14).Value = -40040
$ws.Cells.Item(134, 8).Value = 9436006
$ws.Cells.Item(134, 9).Value = 13890566
$ws.Cells.Item(134, 11).Value = 41671698
$ws.Cells.Item(134, 13).Value = -41669163
$ws.Cells.Item(135, 8).Value = 17096.486
$ws.Cells.Item(135, 10).Value = 17096.486
$ws.Cells.Item(135, 12).Value = 17096.486
$ws.Cells.Item(135, 14).Value = -27236.486
$ws.Cells.Item(138, 8).Value = 18729.846
$ws.Cells.Item(138, 10).Value = 18729.846
$ws.Cells.Item(138, 12).Value = 18729.846
$ws.Cells.Item(138, 14).Value = -29009.846

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 566.3333
$ws.Cells.Item(4, 9).Value = 71
$ws.Cells.Item(4, 10).Value = 690.1667
$ws.Cells.Item(4, 11).Value = 213
$ws.Cells.Item(4, 12).Value = 2070.5001
$ws.Cells.Item(4, 13).Value = -101
$ws.Cells.Item(4, 14).Value = -2294.5001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1903.9231
$ws.Cells.Item(97, 9).Value = 1304.4445
$ws.Cells.Item(97, 11).Value = 1304.4445
$ws.Cells.Item(97, 13).Value = -808.4445000000001
$ws.Cells.Item(124, 8).Value = 31890
$ws.Cells.Item(124, 10).Value = 31890
$ws.Cells.Item(124, 12).Value = 31890
$ws.Cells.Item(124, 14).Value = -41710
$ws.Cells.Item(128, 8).Value = 30000
$ws.Cells.Item(128, 10).Value = 30000
$ws.Cells.Item(128, 12).Value = 30000
$ws.Cells.Item(128, 14).Value = -39960
$ws.Cells.Item(130, 8).Value = 29966.666
$ws.Cells.Item(130, 10).Value = 29966.666
$ws.Cells.Item(130, 12).Value = 29966.666
$ws.Cells.Item(130, 14).Value = -40006.666
$ws.Cells.Item(133, 8).Value = 17330.76
$ws.Cells.Item(133, 10).Value = 17330.76
$ws.Cells.Item(133, 12).Value = 17330.76
$ws.Cells.Item(133, 14).Value = -27450.76
$ws.Cells.Item(135, 8).Value = 24666.666
$ws.Cells.Item(135, 10).Value = 24666.666
$ws.Cells.Item(135, 12).Value = 24666.666
$ws.Cells.Item(135, 14).Value = -34806.666
$ws.Cells.Item(140, 8).Value = 18015.5
$ws.Cells.Item(140, 10).Value = 18015.5
$ws.Cells.Item(140, 12).Value = 18015.5
$ws.Cells.Item(140, 14).Value = -28375.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(125, 8).Value = 29227.445
$ws.Cells.Item(125, 10).Value = 29227.445
$ws.Cells.Item(125, 12).Value = 29227.445
$ws.Cells.Item(125, 14).Value = -39067.445
$ws.Cells.Item(127, 8).Value = 29780.834
$ws.Cells.Item(127, 10).Value = 29780.834
$ws.Cells.Item(127, 12).Value = 29780.834
$ws.Cells.Item(127, 14).Value = -39700.834
$ws.Cells.Item(134, 8).Value = 17955.5
$ws.Cells.Item(134, 10).Value = 17955.5
$ws.Cells.Item(134, 12).Value = 17955.5
$ws.Cells.Item(134, 14).Value = -28095.5
$ws.Cells.Item(139, 8).Value = 30386.182
$ws.Cells.Item(139, 10).Value = 30386.182
$ws.Cells.Item(139, 12).Value = 30386.182
$ws.Cells.Item(139, 14).Value = -40666.182

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2320
$ws.Cells.Item(107, 9).Value = 400
$ws.Cells.Item(107, 10).Value = 10000
$ws.Cells.Item(107, 11).Value = 1200
$ws.Cells.Item(107, 12).Value = 30000
$ws.Cells.Item(107, 13).Value = 720
$ws.Cells.Item(107, 14).Value = -33840
$ws.Cells.Item(125, 8).Value = 29905
$ws.Cells.Item(125, 10).Value = 29905
$ws.Cells.Item(125, 12).Value = 29905
$ws.Cells.Item(125, 14).Value = -39745
$ws.Cells.Item(128, 8).Value = 39980
$ws.Cells.Item(128, 10).Value = 39980
$ws.Cells.Item(128, 12).Value = 39980
$ws.Cells.Item(128, 14).Value = -49940
$ws.Cells.Item(135, 8).Value = 21195.25
$ws.Cells.Item(135, 10).Value = 21195.25
$ws.Cells.Item(135, 12).Value = 21195.25
$ws.Cells.Item(135, 14).Value = -31335.25
$ws.Cells.Item(138, 8).Value = 30507.908
$ws.Cells.Item(138, 10).Value = 30507.908
$ws.Cells.Item(138, 12).Value = 30507.908
$ws.Cells.Item(138, 14).Value = -40787.908
$ws.Cells.Item(139, 8).Value = 18225.19
$ws.Cells.Item(139, 10).Value = 18225.19
$ws.Cells.Item(139, 12).Value = 18225.19
$ws.Cells.Item(139, 14).Value = -28505.19
$ws.Cells.Item(141, 8).Value = 18733.35
$ws.Cells.Item(141, 10).Value = 18733.35
$ws.Cells.Item(141, 12).Value = 18733.35
$ws.Cells.Item(141, 14).Value = -29093.35
